# Insert a new data row before current row 204 (shifting rows 204-300 down to 205-301)
# and populate it with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204; this pushes the existing row 204 (and everything
# below it) down by one row, matching dimension growing from A1:R300 to A1:R301.
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record's values.
$ws.Range("A204").Value = 10
$ws.Range("B204").Value = "Vega Modelo de Temuco"
$ws.Range("C204").Value = "La Araucanía"
$ws.Range("D204").Value = 44518
$ws.Range("E204").Value = 9
$ws.Range("F204").Value = 100112032
$ws.Range("G204").Value = "Zapallo italiano"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 125
$ws.Range("K204").Value = 10000
$ws.Range("L204").Value = 10000
$ws.Range("M204").Value = 10000
$ws.Range("N204").Value = '$/caja 60 unidades'
$ws.Range("O204").Value = "Región del Maule"
$ws.Range("P204").Value = 167
$ws.Range("Q204").Value = 60
$ws.Range("R204").Value = "Hortaliza"
